$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.755.26"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "'1.749.82"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'235.95"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.5055"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").Value = "'40.72"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").Value = "'0.2700"
$ws.Range("E9").Value = "  +13.61%  "
$ws.Range("D10").Value = "'0.06187"
$ws.Range("E10").Value = "  +3.77%  "
$ws.Range("D11").Value = "'1.752.76"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "'0.06929"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "'15.43"
$ws.Range("E13").Value = "  +5.71%  "
$ws.Range("D14").Value = "'0.6098"
$ws.Range("E14").Value = "  +6.52%  "
$ws.Range("D15").Value = "'78.30"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'4.476"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'25.785.49"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'11.65"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").Value = "'0.000006692"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").Value = "'1.978.14"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").Value = "'8.209"
$ws.Range("E24").Value = "  +5.83%  "
$ws.Range("D25").Value = "'5.161"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "'137.27"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "'15.13"
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("D29").Value = "'1.782"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").Value = "'102.26"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").Value = "'0.08270"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "'3.707"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").Value = "'3.394"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'0.04363"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'0.9997"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").Value = "'0.6008"
$ws.Range("D39").Value = "'2.692"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'1.946"
$ws.Range("E40").Value = "  -6.32%  "
$ws.Range("D41").Value = "'0.01549"
$ws.Range("E41").Value = "  +4.72%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'102.11"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "'0.7509"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "'4.846"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("D47").Value = "'0.05498"
$ws.Range("E47").Value = "  +7.63%  "
$ws.Range("D48").Value = "'0.1084"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").Value = "'30.21"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'5.916"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.16%  "
